# Redesign libro diario read algorithm
# Mark the first "credit" (Haber) row of each journal entry block in
# column A with an asterisk "*", so the read algorithm can detect where
# the credit side of an entry begins.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(5, 8, 11, 19, 22, 29, 33, 36)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 1).Value = "*"
}
